$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.685.13"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "3.422.84"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "3.430.08"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "4.012.49"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("D17").Value = "64.621.54"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "3.406.78"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  -4.26%  "
$ws.Range("E27").Value = "  +8.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0762"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "2.878.32"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.772"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "321.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.49%  "
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.07%  "
